$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3 get new ID/NAME/LASTNAME values. The ID column ("56" / "78")
# must be stored as text (not a number), so temporarily force a text format
# on those two cells while entering the values, then clear it again so the
# resulting cells carry no special styling.
$idRange = $ws.Range("A2:A3")
$idRange.NumberFormat = "@"

$ws.Range("A2").Value = "56"
$ws.Range("B2").Value = "hari ram"
$ws.Range("C2").Value = "panchal"

$ws.Range("A3").Value = "78"
$ws.Range("B3").Value = "mohan lal"
$ws.Range("C3").Value = "savita"

$idRange.ClearFormats()

# Row 4: first name "dinesh sengar", last name "amra", last column is an
# empty (but present) text cell. Enter it with a leading apostrophe so it is
# stored as a real (empty) text value instead of being treated as blank,
# then clear the resulting formatting.
$ws.Range("A4").Value = "dinesh sengar"
$ws.Range("B4").Value = "amra"
$c4 = $ws.Range("C4")
$c4.Value = "'"
$c4.ClearFormats()

# The former rows 5 and 6 are removed entirely.
$ws.Rows("5:6").Delete()
